$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Title paragraph: "Xnor.ai acquisition by Apple" -> add " (#2)" run and
#    move the "_GoBack" bookmark here (immediately after the new run).
# ---------------------------------------------------------------------------

# Insert a temporary bookmark right after the existing title text. This acts
# as a hard boundary so the text we insert next stays in its own run instead
# of being coalesced into the existing "Xnor.ai acquisition by Apple" run.
$d.Bookmarks.Add("_TempSplit", $d.Range(28, 28))

# Insert the new run text, plus a throw-away placeholder character. The
# placeholder keeps the upcoming bookmark insertion point away from the
# paragraph's very end (which the engine treats specially and would shove the
# bookmark into the following paragraph).
$d.Range(28, 28).InsertAfter(" (#2)X")

# Re-create the "_GoBack" bookmark immediately after " (#2)" (but before the
# placeholder character).
$d.Bookmarks.Add("_GoBack", $d.Range(33, 33))

# Remove the placeholder character.
$d.Range(33, 34).Delete()

# Remove the temporary splitting bookmark.
$d.Bookmarks("_TempSplit").Delete()

# NOTE: Word only allows a single bookmark per name; re-adding "_GoBack" above
# already relocated (rather than duplicated) the old bookmark that used to
# sit on the "Apparently Xnor's business model..." paragraph, so there is
# nothing further to clean up there.

# ---------------------------------------------------------------------------
# 2. Merge runs that are split apart for no reason (ghost edits where the
#    text content doesn't change but the run boundaries collapse).
# ---------------------------------------------------------------------------

# "Apple have reportedly acquired xnor.ai" + " " + "for " -> " for "
$d.Content.Find.Execute("for ", $true, $false, $false, $false, $false, $true, 1, $false, "for ", 2) | Out-Null

# "The company" + " " + "is" + " well known..." -> "The company is well known..."
$d.Content.Find.Execute("The company is well known", $true, $false, $false, $false, $false, $true, 1, $false, "The company is well known", 2) | Out-Null

# "You " + "require" + " more storage..." -> "You require more storage..."
$d.Content.Find.Execute("You require more storage", $true, $false, $false, $false, $false, $true, 1, $false, "You require more storage", 2) | Out-Null

# " off against increased weight storage and wiring" + " to transmit..." -> combined
$d.Content.Find.Execute(" off against increased weight storage and wiring to transmit", $true, $false, $false, $false, $false, $true, 1, $false, " off against increased weight storage and wiring to transmit", 2) | Out-Null

# " power inference for always-on..." + " are probably the key drivers..." -> combined
$d.Content.Find.Execute(" power inference for always-on features like wake on voice or wake on face are probably", $true, $false, $false, $false, $false, $true, 1, $false, " power inference for always-on features like wake on voice or wake on face are probably", 2) | Out-Null
